$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
$ws.Range("B2").Value = "Thu Jun 19 18:55:41 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:56:25 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:57:06 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:57:43 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
$ws.Range("D2").Value = "Thu Jun 19 18:58:27 IST 2025"
$ws.Range("D3").Value = "Thu Jun 19 18:59:06 IST 2025"
$ws.Range("D4").Value = "Thu Jun 19 18:59:44 IST 2025"
$ws.Range("D5").Value = "Thu Jun 19 19:00:24 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
$ws.Range("B2").Value = "Thu Jun 19 19:01:03 IST 2025"
$ws.Range("B3").Value = "Thu Jun 19 19:01:42 IST 2025"
$ws.Range("B4").Value = "Thu Jun 19 19:02:25 IST 2025"
$ws.Range("B5").Value = "Thu Jun 19 19:03:01 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
$ws.Range("B2").Value = "Thu Jun 19 18:47:53 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:48:32 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:49:14 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:49:54 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item("VT-SaleCredit-NoCF-Generic")
$ws.Range("B2").Value = "Thu Jun 19 18:50:36 IST 2025"
$ws.Range("B3").Value = "Thu Jun 19 18:51:11 IST 2025"
$ws.Range("B4").Value = "Thu Jun 19 18:51:52 IST 2025"
$ws.Range("B5").Value = "Thu Jun 19 18:52:28 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
$ws.Range("B2").Value = "Thu Jun 19 18:53:11 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:53:47 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:54:24 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:55:00 IST 2025"
$ws.Range("C5").Value = "Pass"

$ws = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
$ws.Range("D2").Value = "Mon Jun 16 19:01:35 IST 2025"
$ws.Range("D3").Value = "Mon Jun 16 19:03:33 IST 2025"
$ws.Range("B4").Value = "Tue Jun 17 19:51:35 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("D4").Value = "Thu Jun 19 18:08:57 IST 2025"
$ws.Range("D5").Value = "Mon Jun 16 19:05:48 IST 2025"
$ws.Range("D6").Value = "Mon Jun 16 19:06:58 IST 2025"
$ws.Range("D7").Value = "Mon Jun 16 19:08:16 IST 2025"

$ws = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
$ws.Range("B2").Value = "Tue Jun 17 19:55:48 IST 2025"
$ws.Range("D2").Value = "Thu Jun 19 19:52:35 IST 2025"
$ws.Range("B3").Value = "Tue Jun 17 19:57:04 IST 2025"
$ws.Range("D3").Value = "Thu Jun 19 19:53:30 IST 2025"
$ws.Range("B4").Value = "Tue Jun 17 19:58:12 IST 2025"
$ws.Range("D4").Value = "Thu Jun 19 19:54:20 IST 2025"
$ws.Range("B5").Value = "Tue Jun 17 19:59:15 IST 2025"
$ws.Range("D5").Value = "Thu Jun 19 19:55:10 IST 2025"
$ws.Range("B6").Value = "Tue Jun 17 20:00:28 IST 2025"
$ws.Range("D6").Value = "Thu Jun 19 19:55:59 IST 2025"
$ws.Range("B7").Value = "Tue Jun 17 20:01:32 IST 2025"
$ws.Range("D7").Value = "Thu Jun 19 19:56:43 IST 2025"

$ws = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")
$ws.Range("B2").Value = "Thu Jun 19 18:44:00 IST 2025"
$ws.Range("C2").Value = "Pass"
$ws.Range("B3").Value = "Thu Jun 19 18:44:43 IST 2025"
$ws.Range("C3").Value = "Pass"
$ws.Range("B4").Value = "Thu Jun 19 18:45:22 IST 2025"
$ws.Range("C4").Value = "Pass"
$ws.Range("B5").Value = "Thu Jun 19 18:46:00 IST 2025"
$ws.Range("C5").Value = "Pass"
$ws.Range("B6").Value = "Thu Jun 19 18:46:37 IST 2025"
$ws.Range("C6").Value = "Pass"
$ws.Range("B7").Value = "Thu Jun 19 18:47:17 IST 2025"
$ws.Range("C7").Value = "Pass"
